$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Minimum module size" values for the rows that remain (network 1-8, GWAS model).
$ws.Range("D2").Value = 10
$ws.Range("D3").Value = 10
$ws.Range("D4").Value = 25
$ws.Range("D5").Value = 25
$ws.Range("D6").Value = 10
$ws.Range("D7").Value = 10
$ws.Range("D8").Value = 25
$ws.Range("D9").Value = 25

# 2. Remove the RR-BLUP block (old rows 10-21), shrinking the table down to network 1-8.
$ws.Range("A10:E21").EntireRow.Delete()

# 3. Give the new last row (row 9) the same bottom-border / centered look as the
#    header row, but without bold text (matches the existing "last row" style
#    used previously by row 21).
$ws.Range("A9:E9").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)
$ws.Range("A9:E9").Font.Bold = $false
$ws.Range("A1:E1").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)
$ws.Range("A9:E9").Font.Bold = $false

# 4. Add the new "Markers in network" column.
$ws.Range("A1:E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Markers in network"

$ws.Range("F2").Value = 2020
$ws.Range("F3").Value = 2020
$ws.Range("F4").Value = 2020
$ws.Range("F5").Value = 2020
$ws.Range("F6").Value = 1832
$ws.Range("F7").Value = 1832
$ws.Range("F8").Value = 1832
$ws.Range("F9").Value = 1832

$ws.Range("A2:E8").Copy()
$ws.Range("A2:E8").PasteSpecial(-4122)

$ws.Columns.Item(6).ColumnWidth = 17.5

# 5. Re-apply the bottom border/centered (non-bold) look to row 9 after the
#    column F value landed, and make sure F9 matches it too.
$ws.Range("A9:F9").Copy()
$ws.Range("A9:F9").PasteSpecial(-4122)
$ws.Range("A9:F9").Font.Bold = $false

# 6. Update the selection to match the saved workbook state.
$ws.Range("G8").Select()
